$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: update the CDC reference date (B8) ---
# Formula in C8 ("shared" formula si=1) recalculates automatically.
$ws.Range("B8").Value = 44659

# --- Row 10: new "IPS" establishment row -----------------------------------
# Added first so the new shared-string / fill entries are allocated in the
# same order as in the target workbook (IPS before MILL).
$ws.Range("A10").Value = "IPS"
$ws.Range("A10").Interior.Color = 15773696   # RGB(0,176,240) -> FF00B0F0
$ws.Range("B10").Value = 44592

# --- Row 9: new "MILL" establishment row ------------------------------------
$ws.Range("A9").Value = "MILL"
$ws.Range("A9").Interior.Color = 12611584    # RGB(0,112,192) -> FF0070C0
$ws.Range("B9").Value = 44656

# --- Selection / active cell moved to C8 ------------------------------------
$ws.Range("C8").Select()
